$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2 through 498) all currently hold the date-serial value
# 45180; update them to 45181 to reflect the new "Förändrad" date.
$ws.Range("C2:C498").Value = 45181
